# Update brochure file link addresses: "offshore/seapae/reach" -> "offshore/region-1/reach"
# (commit message: "update brochure file link address")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New target URLs (region-1 instead of seapae)
$tuitionUrl    = "https://github.com/Viskee-Consultancy/Viskee-Consultancy-Configuration/raw/master/brochures/pdf/offshore/region-1/reach/Reach_Tuition_Course_Fees_2021_v1.0.pdf"
$promotionsUrl = "https://github.com/Viskee-Consultancy/Viskee-Consultancy-Configuration/raw/master/brochures/pdf/offshore/region-1/reach/Reach-Offshore-SEAPAE-Q4-Promotions-1OCT-31DEC21_VOL-1.1.pdf"

# Update the cell text to point at the new location
$ws.Range("B2").Value = $tuitionUrl
$ws.Range("B3").Value = $promotionsUrl

# Turn those cells into real hyperlinks pointing at the (updated) URL text.
# This also creates the "Hyperlink" cell style (font/xf/cellStyle) and applies
# it to B2/B3, matching the styles.xml additions in the diff.
$ws.Hyperlinks.Add($ws.Range("B2"), $ws.Range("B2").Value2)
$ws.Hyperlinks.Add($ws.Range("B3"), $ws.Range("B3").Value2)

# Widen column B to fit the (now longer) region-1 link text.
$ws.Columns.Item(2).ColumnWidth = 179.16071428571428
